$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''26.763.37'
$ws.Range("E2").Value = '  +1.28%  '
$ws.Range("D3").Value = '''1.725.84'
$ws.Range("E3").Value = '  +0.05%  '
$ws.Range("D4").Value = '''0.9974'
$ws.Range("E4").Value = '  -0.21%  '
$ws.Range("D5").Value = '''240.57'
$ws.Range("E5").Value = '  -1.16%  '
$ws.Range("D6").Value = '''0.9978'
$ws.Range("E6").Value = '  -0.22%  '
$ws.Range("D7").Value = '''0.4859'
$ws.Range("E7").Value = '  -1.07%  '
$ws.Range("D8").Value = '''0.2584'
$ws.Range("E8").Value = '  -0.91%  '
$ws.Range("E9").Value = '  -0.12%  '
$ws.Range("D10").Value = '''1.728.02'
$ws.Range("E10").Value = '  +0.23%  '
$ws.Range("E11").Value = '  +2.99%  '
$ws.Range("D12").Value = '''0.06870'
$ws.Range("E12").Value = '  -1.68%  '
$ws.Range("D13").Value = '''0.6077'
$ws.Range("E13").Value = '  +1.34%  '
$ws.Range("D14").Value = '''4.479'
$ws.Range("E14").Value = '  -1.43%  '
$ws.Range("D15").Value = '''76.88'
$ws.Range("E15").Value = '  -0.73%  '
$ws.Range("D16").Value = '''0.9979'
$ws.Range("E16").Value = '  -0.21%  '
$ws.Range("D17").Value = '''26.564.37'
$ws.Range("E17").Value = '  +0.55%  '
$ws.Range("D18").Value = '''0.9974'
$ws.Range("E18").Value = '  -0.23%  '
$ws.Range("D19").Value = '''0.000007150'
$ws.Range("E19").Value = '  -1.00%  '
$ws.Range("D20").Value = '''11.44'
$ws.Range("E20").Value = '  +0.76%  '
$ws.Range("D21").Value = '''1.949.99'
$ws.Range("E21").Value = '  +0.29%  '
$ws.Range("D22").Value = '''4.430'
$ws.Range("E22").Value = '  -0.72%  '
$ws.Range("D23").Value = '''8.578'
$ws.Range("E23").Value = '  -0.17%  '
$ws.Range("D24").Value = '''5.090'
$ws.Range("E24").Value = '  -1.36%  '
$ws.Range("D25").Value = '''137.50'
$ws.Range("E25").Value = '  -0.06%  '
$ws.Range("D26").Value = '''15.24'
$ws.Range("E26").Value = '  -0.21%  '
$ws.Range("E27").Value = '  +2.96%  '
$ws.Range("D28").Value = '''105.96'
$ws.Range("E28").Value = '  -0.92%  '
$ws.Range("E29").Value = '  -1.89%  '
$ws.Range("D30").Value = '''3.998'
$ws.Range("D31").Value = '''0.07930'
$ws.Range("E31").Value = '  -0.90%  '
$ws.Range("E32").Value = '  +0.12%  '
$ws.Range("E33").Value = '  -0.67%  '
$ws.Range("B34").Value = 'HuobiToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D34").Value = '''2.596'
$ws.Range("E34").Value = '  -0.26%  '
$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D35").Value = '''1.006'
$ws.Range("E35").Value = '  +0.24%  '
$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").Value = '''0.6211'
$ws.Range("E36").Value = '  -0.79%  '
$ws.Range("B37").Value = 'TrustWalletToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D37").Value = '''0.9229'
$ws.Range("E37").Value = '  -2.09%  '
$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D38").Value = '''2.026'
$ws.Range("E38").Value = '  +3.88%  '
$ws.Range("B39").Value = 'MXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D39").Value = '''2.443'
$ws.Range("E39").Value = '  +2.25%  '
$ws.Range("B40").Value = 'PaxDollar'
$ws.Range("C40").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D40").Value = '''0.9972'
$ws.Range("E40").Value = '  -0.24%  '
$ws.Range("B41").Value = 'VeChain'
$ws.Range("C41").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D41").Value = '''0.01495'
$ws.Range("E41").Value = '  +0.77%  '
$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").Value = '''5.652'
$ws.Range("E42").Value = '  +6.29%  '
$ws.Range("B43").Value = 'Quant'
$ws.Range("C43").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D43").Value = '''99.84'
$ws.Range("E43").Value = '  +0.22%  '
$ws.Range("B44").Value = 'TheSandbox'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D44").Value = '''0.3841'
$ws.Range("E44").Value = '  -0.42%  '
$ws.Range("B45").Value = 'Aptos'
$ws.Range("C45").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D45").Value = '''6.859'
$ws.Range("E45").Value = '  +0.36%  '
$ws.Range("B46").Value = 'Algorand'
$ws.Range("C46").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D46").Value = '''0.1156'
$ws.Range("E46").Value = '  -1.35%  '
$ws.Range("B47").Value = 'Cronos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D47").Value = '''0.05380'
$ws.Range("E47").Value = '  +0.29%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").Value = '''7.861'
$ws.Range("E48").Value = '  +1.46%  '
$ws.Range("B49").Value = 'Elrond'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D49").Value = '''30.10'
$ws.Range("E49").Value = '  -0.43%  '
$ws.Range("B50").Value = 'NEARProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D50").Value = '''1.234'
$ws.Range("E50").Value = '  -0.19%  '
$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D51").Value = '''51.49'
$ws.Range("E51").Value = '  +1.18%  '
